$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1.2 - add a title row above the data table
$ws.Range("B1").Value = "Profile"
$ws.Range("C1").Value = "Avg Scores"
$ws.Range("D1").Value = "Students Number"
$ws.Range("E1").Value = "University Number"
$ws.Range("F1").Value = "University List"

# The "Avg Scores" column now holds text (including "NaN" for profiles
# whose average could not be computed) instead of numbers/#NUM! errors.
# Force the numeric-looking values to be stored as text, then drop the
# temporary text format so the cells keep the default style.
$ws.Range("C2:C5").NumberFormat = "@"
$ws.Range("C2").Value = "4.537499964237213"
$ws.Range("C3").Value = "4.333333492279053"
$ws.Range("C4").Value = "NaN"
$ws.Range("C5").Value = "NaN"
$ws.Range("C2:C5").ClearFormats()
